# ---------------------------------------------------------------------------
# Lagt til prosjektfiler fra BB, og begynt å fylle inn systemdata - basecase
#
# 1. Renames the original sheet "Sheet1" -> "Avrundet" (rounded values,
#    already present).
# 2. Adds a new sheet "Uten avrunding" (unrounded values) right after it,
#    with the same layout/headers but the raw (unrounded) computed numbers,
#    and with the "Shunt Admittans" columns relabelled / reordered slightly.
# 3. Updates the selection on both sheets and makes the new sheet active.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Avrundet"

# New sheet, inserted right after "Avrundet"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Uten avrunding"

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$header = @(
    "Linje",
    "Impedans (ohm)",
    "Impedans (p.u.)",
    "Admittans (p.u.)",
    "Kapasitans (nF)",
    "Shunt Impedans (ohm)",
    "Shunt Admittans (ohm)",
    "Shunt Admittans (p.u.)",
    "Shunt Admittans (p.u.) half"
)
for ($c = 1; $c -le 9; $c++) { $ws2.Cells.Item(1, $c).Value = $header[$c - 1] }

# ---------------------------------------------------------------------------
# Data rows (unrounded line-impedance calculations)
# ---------------------------------------------------------------------------
$rows = @(
    @("1-2", "(0.9159999999999999+9.847j)", "(0.010177777777777777+0.10941111111111111j)", "(0.8429235398105763-9.061428052963697j)", 199.23, "-15977.005781448112j", "6.258995043746944e-05j", "0.005633095539372249j", "0.0028165477696861247j"),
    @("2-3", "(0.752+8.084j)", "(0.008355555555555555+0.08982222222222222j)", "(1.0267526096628827-11.03759055387599j)", 163.56, "-19461.352786976684j", "5.1383889442114665e-05j", "0.00462455004979032j", "0.00231227502489516j"),
    @("3-4", "(1.912+20.554j)", "(0.021244444444444444+0.22837777777777776j)", "(0.4038273862272426-4.341144401942858j)", 415.86, "-7654.255907848572j", "0.00013064627209218512j", "0.011758164488296662j", "0.005879082244148331j"),
    @("4-5", "(0.363+3.993j)", "(0.004033333333333333+0.044366666666666665j)", "(2.0322449532583664-22.35469448584203j)", 136.73, "-23280.17890615012j", "4.2954996352533256e-05j", "0.003865949671727993j", "0.0019329748358639966j"),
    @("5-6", "(1.4240000000000002+22.784000000000002j)", "(0.015822222222222224+0.2531555555555556j)", "(0.24592314082105537-3.934770253136886j)", 811.68, "-3921.6179551521614j", "0.0002549967925065764j", "0.022949711325591874j", "0.011474855662795937j"),
    @("6-7", "(1.088+17.408j)", "(0.01208888888888889+0.19342222222222225j)", "(0.3218699931334401-5.149919890135042j)", 620.16, "-5132.7058530668j", "0.0001948290100050246j", "0.017534610900452215j", "0.008767305450226108j"),
    @("7-8", "(1.976+21.241999999999997j)", "(0.021955555555555555+0.2360222222222222j)", "(0.3907479567138098-4.200540534673455j)", 429.78, "-7406.34478532716j", "0.00013501936906598212j", "0.012151743215938392j", "0.006075871607969196j"),
    @("1-8", "(4.232+45.494j)", "(0.04702222222222223+0.5054888888888889j)", "(0.182447533664104-1.9613109868891176j)", 920.46, "-3458.160986721755j", "0.0002891710373923261j", "0.02602539336530935j", "0.013012696682654675j"),
    @("1-6", "(4.668+50.181j)", "(0.05186666666666667+0.5575666666666667j)", "(0.16540659007422623-1.7781208432979319j)", 1015.29, "-3135.162231321008j", "0.00031896276052631814j", "0.028706648447368633j", "0.014353324223684316j")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 1; $c -le 9; $c++) {
        $ws2.Cells.Item($r + 2, $c).Value = $rowVals[$c - 1]
    }
}

# ---------------------------------------------------------------------------
# Formatting: copy header style + column widths from "Avrundet"
# ---------------------------------------------------------------------------
$ws1.Range("A1:I1").Copy()
$ws2.Range("A1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Columns.Item(2).ColumnWidth = 41.8776041666667
$ws2.Columns.Item(3).ColumnWidth = 44.0221354166667
$ws2.Columns.Item(4).ColumnWidth = 39.5924479166667
$ws2.Columns.Item(5).ColumnWidth = 13.0221354166667
$ws2.Columns.Item(6).ColumnWidth = 21.8776041666667
$ws2.Columns.Item(7).ColumnWidth = 27.0221354166667
$ws2.Columns.Item(8).ColumnWidth = 22.5924479166667
$ws2.Columns.Item(9).ColumnWidth = 24.3072916666667

# ---------------------------------------------------------------------------
# Page margins on the new sheet
# ---------------------------------------------------------------------------
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# Selections: "Avrundet" keeps a selection on I12 (no longer the active tab);
# "Uten avrunding" becomes the active tab with its selection on F18.
# ---------------------------------------------------------------------------
[void]$ws1.Range("I12").Select()
$ws2.Activate()
[void]$ws2.Range("F18").Select()
